# houseparts.xlsx: merge in "Materials and suppliers" sheet/table, rename
# the original sheet, and refresh the House parts table (External Walls
# replaces Walls, Wood replaces Tiles, Glass replaces Pavers).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# --- Rename the original sheet ---------------------------------------
$ws1.Name = "House parts and materials"

# --- Update the House-parts table data --------------------------------
# Before:                              After:
#  Roof      | Shingles | Metal | Tiles    Roof           | Shingles | Metal | Wood
#  Walls     | Brick    | Concrete | Wood  External Walls | Brick    | Concrete | Wood
#  Foundation| Concrete | Stone | Pavers   Foundation     | Concrete | Stone | Glass
$ws1.Range("D2").Value2 = "Wood"
$ws1.Range("A3").Value2 = "External Walls"
$ws1.Range("D4").Value2 = "Glass"

# --- Add the new "Materials and suppliers" sheet, right after sheet 1 -
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Materials and suppliers"

# --- Populate the Materials/Suppliers/Cost-per-unit table -------------
$ws2.Range("C3").Value2 = "Material"
$ws2.Range("D3").Value2 = "Supplier"
$ws2.Range("E3").Value2 = "Cost per unit"

$ws2.Range("C4").Value2 = "Shingles"
$ws2.Range("D4").Value2 = "ABC Shingles"
$ws2.Range("E4").Value2 = 10

$ws2.Range("C5").Value2 = "Shingles"
$ws2.Range("D5").Value2 = "DEF Shingles"
$ws2.Range("E5").Value2 = 20

$ws2.Range("C6").Value2 = "Concrete"
$ws2.Range("D6").Value2 = "ABC Concrete"
$ws2.Range("E6").Value2 = 40

$ws2.Range("C7").Value2 = "Concrete"
$ws2.Range("D7").Value2 = "DEF Concrete"
$ws2.Range("E7").Value2 = 34

# --- Turn that range into a table, styled like the rest of the workbook
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("C3:E7"), $null, 1)
$tbl2.Name = "Table2"
$tbl2.TableStyle = "TableStyleLight2"

# --- Restore selections / active sheet to match the saved workbook ----
$ws2.Range("E8").Select()
$ws1.Range("D3").Select()
$ws1.Activate()
